$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "OK" markers for rows whose todo item has been done (col C)
$ws.Range("C6").Value = "OK"
$ws.Range("C7").Value = "OK"
$ws.Range("C8").Value = "OK"
$ws.Range("C10").Value = "OK"

# New column D notes (order matters for shared-string table indices)
$ws.Range("D9").Value = "?"
$ws.Range("D2").Value = "bug : taille d'image "

# Update selection to match the recorded cursor position
$ws.Range("D3").Select()
